$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")

# --- Row 48: fill in the remaining columns for the existing TestCaseID ---
$ws.Range("B48").Value = "Yes"
$ws.Range("C48").Value = "Android"
$ws.Range("D48").Value = "user2046417"
$ws.Range("E48").Value = "Kony@1234"
[void]$ws.Hyperlinks.Add($ws.Range("E48"), "mailto:Kony@1234")

# --- Row 49: new test case - Stop Payment via Services ---
$ws.Range("A49").Value = "C23629_VerifyStopPaymentViaServices"
$ws.Range("B49").Value = "Yes"
$ws.Range("C49").Value = "Android"
$ws.Range("D49").Value = "user2046417"
$ws.Range("E49").Value = "Kony@1234"
[void]$ws.Hyperlinks.Add($ws.Range("E49"), "mailto:Kony@1234")

# --- Row 50: new test case - TDECU Rates via Services ---
$ws.Range("A50").Value = "C23635_VerifyTDECURatesViaServices"
$ws.Range("B50").Value = "Yes"
$ws.Range("C50").Value = "Android"
$ws.Range("D50").Value = "user2046417"
$ws.Range("E50").Value = "Kony@1234"
[void]$ws.Hyperlinks.Add($ws.Range("E50"), "mailto:Kony@1234")

# --- Row 51: new test case - Credit Card Information via Services ---
$ws.Range("A51").Value = "C23634_VerifyCreditCardInformationViaServices"
$ws.Range("B51").Value = "Yes"
$ws.Range("C51").Value = "Android"
$ws.Range("D51").Value = "user2046417"
$ws.Range("E51").Value = "Kony@1234"
[void]$ws.Hyperlinks.Add($ws.Range("E51"), "mailto:Kony@1234")

# Re-apply the shared "Hyperlink" cell style so the newly linked cells reuse the
# workbook's existing style record instead of a freshly minted duplicate.
$ws.Range("E48:E51").Style = "Hyperlink"

# Update the selection / active cell to the last new row, matching the author's
# on-screen state when the change was saved.
$ws.Activate()
[void]$ws.Range("B51:F51").Select()

# The "General" sheet becomes the active tab of the workbook (it was "Data" before).
